# The example workbook previously contained random numbers in A1:B2.
# Replace that with a single text cell (A1 = "This is neat!") and remove
# the other three cells (B1, A2, B2) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old numeric cells that are no longer part of the sheet.
$ws.Range("B1:B2").ClearContents()
$ws.Range("A2").ClearContents()

# A1 becomes a shared string with the new text.
$ws.Range("A1").Value = "This is neat!"
